$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Existing-cell edits per diff ---
$ws.Range("Q55").Value = 0
$ws.Range("Q63").Value = 0
$ws.Range("Q66").Value = 0
$ws.Range("O1121").Value = 3
$ws.Range("R1123").Value = 0
$ws.Range("R1124").Value = 0

# --- New rows 1125-1140 ---
$ws.Range("A1125").Value = 45474
$ws.Range("A1125").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B1125").Value = 119.3499984741211
$ws.Range("C1125").Value = 120.1900024414062
$ws.Range("D1125").Value = 115.5
$ws.Range("E1125").Value = 117.7600021362305
$ws.Range("F1125").Value = 117.7600021362305
$ws.Range("G1125").Value = 100131698
$ws.Range("H1125").Value = 2024
$ws.Range("I1125").Value = 7
$ws.Range("J1125").Value = 1
$ws.Range("K1125").Value = 0
$ws.Range("L1125").Value = 0
$ws.Range("M1125").Value = 0
$ws.Range("N1125").Value = 27
$ws.Range("O1125").Value = 0
$ws.Range("P1125").Value = 0
$ws.Range("Q1125").Value = 0

$ws.Range("A1126").Value = 45481
$ws.Range("A1126").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B1126").Value = 117.5
$ws.Range("C1126").Value = 118.0800018310547
$ws.Range("D1126").Value = 112.3000030517578
$ws.Range("E1126").Value = 112.7200012207031
$ws.Range("F1126").Value = 112.7200012207031
$ws.Range("G1126").Value = 125834974
$ws.Range("H1126").Value = 2024
$ws.Range("I1126").Value = 7
$ws.Range("J1126").Value = 8
$ws.Range("K1126").Value = 0
$ws.Range("L1126").Value = 0
$ws.Range("M1126").Value = 0
$ws.Range("N1126").Value = 28
$ws.Range("O1126").Value = 0
$ws.Range("P1126").Value = 0
$ws.Range("Q1126").Value = 0

$ws.Range("A1127").Value = 45488
$ws.Range("A1127").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B1127").Value = 113.25
$ws.Range("C1127").Value = 118.1100006103516
$ws.Range("D1127").Value = 112.3300018310547
$ws.Range("E1127").Value = 112.8899993896484
$ws.Range("F1127").Value = 112.8899993896484
$ws.Range("G1127").Value = 141782015
$ws.Range("H1127").Value = 2024
$ws.Range("I1127").Value = 7
$ws.Range("J1127").Value = 15
$ws.Range("K1127").Value = 0
$ws.Range("L1127").Value = 0
$ws.Range("M1127").Value = 0
$ws.Range("N1127").Value = 29
$ws.Range("O1127").Value = 0
$ws.Range("P1127").Value = 0
$ws.Range("Q1127").Value = 0

$ws.Range("A1128").Value = 45495
$ws.Range("A1128").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B1128").Value = 112.8000030517578
$ws.Range("C1128").Value = 115.4899978637695
$ws.Range("D1128").Value = 109
$ws.Range("E1128").Value = 113.8600006103516
$ws.Range("F1128").Value = 113.8600006103516
$ws.Range("G1128").Value = 203853368
$ws.Range("H1128").Value = 2024
$ws.Range("I1128").Value = 7
$ws.Range("J1128").Value = 22
$ws.Range("K1128").Value = 0
$ws.Range("L1128").Value = 0
$ws.Range("M1128").Value = 0
$ws.Range("N1128").Value = 30
$ws.Range("O1128").Value = 0
$ws.Range("P1128").Value = 0
$ws.Range("Q1128").Value = 0

$ws.Range("A1129").Value = 45502
$ws.Range("A1129").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B1129").Value = 115.1999969482422
$ws.Range("C1129").Value = 117.8199996948242
$ws.Range("D1129").Value = 110.0500030517578
$ws.Range("E1129").Value = 110.4899978637695
$ws.Range("F1129").Value = 110.4899978637695
$ws.Range("G1129").Value = 139329247
$ws.Range("H1129").Value = 2024
$ws.Range("I1129").Value = 7
$ws.Range("J1129").Value = 29
$ws.Range("K1129").Value = 0
$ws.Range("L1129").Value = 0
$ws.Range("M1129").Value = 0
$ws.Range("N1129").Value = 31
$ws.Range("O1129").Value = 0
$ws.Range("P1129").Value = 0
$ws.Range("Q1129").Value = 0

$ws.Range("A1130").Value = 45509
$ws.Range("A1130").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B1130").Value = 108.4400024414062
$ws.Range("C1130").Value = 112.5400009155273
$ws.Range("D1130").Value = 104.0400009155273
$ws.Range("E1130").Value = 110.6500015258789
$ws.Range("F1130").Value = 110.6500015258789
$ws.Range("G1130").Value = 147956939
$ws.Range("H1130").Value = 2024
$ws.Range("I1130").Value = 8
$ws.Range("J1130").Value = 5
$ws.Range("K1130").Value = 0
$ws.Range("L1130").Value = 0
$ws.Range("M1130").Value = 0
$ws.Range("N1130").Value = 32
$ws.Range("O1130").Value = 0
$ws.Range("P1130").Value = 0
$ws.Range("Q1130").Value = 0

$ws.Range("A1131").Value = 45516
$ws.Range("A1131").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B1131").Value = 109.9000015258789
$ws.Range("C1131").Value = 110.7600021362305
$ws.Range("D1131").Value = 105.5
$ws.Range("E1131").Value = 107.620002746582
$ws.Range("F1131").Value = 107.620002746582
$ws.Range("G1131").Value = 85508591
$ws.Range("H1131").Value = 2024
$ws.Range("I1131").Value = 8
$ws.Range("J1131").Value = 12
$ws.Range("K1131").Value = 0
$ws.Range("L1131").Value = 0
$ws.Range("M1131").Value = 0
$ws.Range("N1131").Value = 33
$ws.Range("O1131").Value = 0
$ws.Range("P1131").Value = 0
$ws.Range("Q1131").Value = 0

$ws.Range("A1132").Value = 45523
$ws.Range("A1132").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B1132").Value = 108.3000030517578
$ws.Range("C1132").Value = 112.9899978637695
$ws.Range("D1132").Value = 108.0800018310547
$ws.Range("E1132").Value = 112.0999984741211
$ws.Range("F1132").Value = 112.0999984741211
$ws.Range("G1132").Value = 94610840
$ws.Range("H1132").Value = 2024
$ws.Range("I1132").Value = 8
$ws.Range("J1132").Value = 19
$ws.Range("K1132").Value = 0
$ws.Range("L1132").Value = 0
$ws.Range("M1132").Value = 0
$ws.Range("N1132").Value = 34
$ws.Range("O1132").Value = 0
$ws.Range("P1132").Value = 0
$ws.Range("Q1132").Value = 0

$ws.Range("A1133").Value = 45530
$ws.Range("A1133").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B1133").Value = 113.1600036621094
$ws.Range("C1133").Value = 113.1900024414062
$ws.Range("D1133").Value = 108.8499984741211
$ws.Range("E1133").Value = 111.5299987792969
$ws.Range("F1133").Value = 111.5299987792969
$ws.Range("G1133").Value = 74014250
$ws.Range("H1133").Value = 2024
$ws.Range("I1133").Value = 8
$ws.Range("J1133").Value = 26
$ws.Range("K1133").Value = 0
$ws.Range("L1133").Value = 0
$ws.Range("M1133").Value = 0
$ws.Range("N1133").Value = 35
$ws.Range("O1133").Value = 0
$ws.Range("P1133").Value = 0
$ws.Range("Q1133").Value = 0

$ws.Range("A1134").Value = 45537
$ws.Range("A1134").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B1134").Value = 111.8600006103516
$ws.Range("C1134").Value = 113.4599990844727
$ws.Range("D1134").Value = 102.9499969482422
$ws.Range("E1134").Value = 103.379997253418
$ws.Range("F1134").Value = 103.379997253418
$ws.Range("G1134").Value = 162619891
$ws.Range("H1134").Value = 2024
$ws.Range("I1134").Value = 9
$ws.Range("J1134").Value = 2
$ws.Range("K1134").Value = 0
$ws.Range("L1134").Value = 0
$ws.Range("M1134").Value = 0
$ws.Range("N1134").Value = 36
$ws.Range("O1134").Value = 0
$ws.Range("P1134").Value = 0
$ws.Range("Q1134").Value = 0

$ws.Range("A1135").Value = 45544
$ws.Range("A1135").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B1135").Value = 103.5999984741211
$ws.Range("C1135").Value = 107.25
$ws.Range("D1135").Value = 100.6900024414062
$ws.Range("E1135").Value = 106.3899993896484
$ws.Range("F1135").Value = 106.3899993896484
$ws.Range("G1135").Value = 142400116
$ws.Range("H1135").Value = 2024
$ws.Range("I1135").Value = 9
$ws.Range("J1135").Value = 9
$ws.Range("K1135").Value = 0
$ws.Range("L1135").Value = 0
$ws.Range("M1135").Value = 0
$ws.Range("N1135").Value = 37
$ws.Range("O1135").Value = 2
$ws.Range("P1135").Value = 0
$ws.Range("Q1135").Value = 0

$ws.Range("A1136").Value = 45551
$ws.Range("A1136").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B1136").Value = 106.9000015258789
$ws.Range("C1136").Value = 108.1800003051758
$ws.Range("D1136").Value = 102.7200012207031
$ws.Range("E1136").Value = 104.9599990844727
$ws.Range("F1136").Value = 104.9599990844727
$ws.Range("G1136").Value = 86000137
$ws.Range("H1136").Value = 2024
$ws.Range("I1136").Value = 9
$ws.Range("J1136").Value = 16
$ws.Range("K1136").Value = 0
$ws.Range("L1136").Value = 0
$ws.Range("M1136").Value = 0
$ws.Range("N1136").Value = 38
$ws.Range("O1136").Value = 0
$ws.Range("P1136").Value = 0
$ws.Range("Q1136").Value = 0

$ws.Range("A1137").Value = 45558
$ws.Range("A1137").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B1137").Value = 105.8000030517578
$ws.Range("C1137").Value = 113.3000030517578
$ws.Range("D1137").Value = 105.0500030517578
$ws.Range("E1137").Value = 113.0999984741211
$ws.Range("F1137").Value = 113.0999984741211
$ws.Range("G1137").Value = 240494416
$ws.Range("H1137").Value = 2024
$ws.Range("I1137").Value = 9
$ws.Range("J1137").Value = 23
$ws.Range("K1137").Value = 0
$ws.Range("L1137").Value = 0
$ws.Range("M1137").Value = 0
$ws.Range("N1137").Value = 39
$ws.Range("O1137").Value = 0
$ws.Range("P1137").Value = 0
$ws.Range("Q1137").Value = 0

$ws.Range("A1138").Value = 45565
$ws.Range("A1138").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B1138").Value = 113
$ws.Range("C1138").Value = 113.5
$ws.Range("D1138").Value = 106.5100021362305
$ws.Range("E1138").Value = 107.620002746582
$ws.Range("F1138").Value = 107.620002746582
$ws.Range("G1138").Value = 87128675
$ws.Range("H1138").Value = 2024
$ws.Range("I1138").Value = 9
$ws.Range("J1138").Value = 30
$ws.Range("K1138").Value = 0
$ws.Range("L1138").Value = 0
$ws.Range("M1138").Value = 0
$ws.Range("N1138").Value = 40
$ws.Range("O1138").Value = 0
$ws.Range("P1138").Value = 0
$ws.Range("Q1138").Value = 0

$ws.Range("A1139").Value = 45572
$ws.Range("A1139").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B1139").Value = 108.2399978637695
$ws.Range("C1139").Value = 109.0500030517578
$ws.Range("D1139").Value = 102.3399963378906
$ws.Range("E1139").Value = 104.0599975585938
$ws.Range("F1139").Value = 104.0599975585938
$ws.Range("G1139").Value = 108484220
$ws.Range("H1139").Value = 2024
$ws.Range("I1139").Value = 10
$ws.Range("J1139").Value = 7
$ws.Range("K1139").Value = 0
$ws.Range("L1139").Value = 0
$ws.Range("M1139").Value = 0
$ws.Range("N1139").Value = 41
$ws.Range("O1139").Value = 0
$ws.Range("P1139").Value = 0
$ws.Range("Q1139").Value = 0

$ws.Range("A1140").Value = 45579
$ws.Range("A1140").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B1140").Value = 104.5
$ws.Range("C1140").Value = 105.4400024414062
$ws.Range("D1140").Value = 100.7699966430664
$ws.Range("E1140").Value = 104.6699981689453
$ws.Range("F1140").Value = 104.6699981689453
$ws.Range("G1140").Value = 65738218
$ws.Range("H1140").Value = 2024
$ws.Range("I1140").Value = 10
$ws.Range("J1140").Value = 14
$ws.Range("K1140").Value = 0
$ws.Range("L1140").Value = 0
$ws.Range("M1140").Value = 0
$ws.Range("N1140").Value = 42
$ws.Range("O1140").Value = 0
$ws.Range("P1140").Value = 0
$ws.Range("Q1140").Value = 0
